# Trade #54 closed at 2026-02-17 08:42:28 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: update aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.51   # Current Capital
$summary.Range("B4").Value = -0.48     # Total P&L $
$summary.Range("B5").Value = -0.18     # Total P&L %
$summary.Range("B6").Value = 54        # Total Trades
$summary.Range("B7").Value = 20        # Winning Trades
$summary.Range("B9").Value = 37.04     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.51000000000001   # Capital
$status.Range("D4").Value = 54                  # Trades
$status.Range("E4").Value = -0.48               # P&L $
$status.Range("F4").Value = -0.49               # P&L %
$status.Range("G4").Value = 37.04               # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade (#54) to both "All Trades" and "MarketMaking"
# sheets as a new row 55.
# ---------------------------------------------------------------------
$newRow = @(54, "2026-02-17", "08:42:22", "MarketMaking", "DOWN", 0.36, 0.43, "CLOSED", 19.4444, 0.07000000000000001, 99.51000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 55
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($row, $col).Value = $newRow[$i]
    }
    # The Date column (B) looks like a date literal, so Excel auto-converts
    # it to a date serial number. Force it back to plain text (matching the
    # rest of the column) without leaving a lingering text number-format.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newRow[1]
    $dateCell.ClearFormats()
}
